$d = $word.ActiveDocument

# Locate the paragraph containing "LOB1004: Cálculo II (Requisito fraco)" —
# the three paragraphs that follow it (an empty paragraph, the
# "Ver no Jupiter Salvar em pdf Salvar em docx" paragraph, and the
# "© 2020 . Contact: ..." paragraph) are removed, per the diff.

$rng = $d.Content
$rng.Find.Execute("LOB1004: Cálculo II (Requisito fraco)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$anchor = $rng.Paragraphs.Item(1)
$idx = $anchor.Index

# Delete from the last of the three target paragraphs back to the first, so
# indices of paragraphs we still need to remove stay valid as we go.
$d.Paragraphs.Item($idx + 3).Range.Delete()
$d.Paragraphs.Item($idx + 2).Range.Delete()
$d.Paragraphs.Item($idx + 1).Range.Delete()
